$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6:D6").Copy()
$ws.Range("A7:D7").PasteSpecial(-4104)

$ws.Range("A7").Value = "As a developer, I want to create a prototype of the UI."
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 3
$ws.Range("D7").Value = "The prototype is created to make it easier to implement in Android Studio."

$ws.Range("J5").Select()
